# Update "想去人数" (want-to-go count) values in column F across all four
# sheets, matching the "Update gh-pages to output generated at 456a3b4" diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2467
$ws.Range("F5").Value = 1769
$ws.Range("F6").Value = 110
$ws.Range("F9").Value = 3581
$ws.Range("F10").Value = 1201
$ws.Range("F11").Value = 1586
$ws.Range("F14").Value = 21
$ws.Range("F15").Value = 1374
$ws.Range("F19").Value = 469
$ws.Range("F21").Value = 36
$ws.Range("F22").Value = 5
$ws.Range("F24").Value = 2324
$ws.Range("F25").Value = 159
$ws.Range("F31").Value = 147
$ws.Range("F34").Value = 940

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 29
$ws.Range("F8").Value = 20
$ws.Range("F14").Value = 41
$ws.Range("F22").Value = 141
$ws.Range("F24").Value = 194
$ws.Range("F35").Value = 446
$ws.Range("F39").Value = 30

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 9584
$ws.Range("F10").Value = 3011
$ws.Range("F11").Value = 524
$ws.Range("F13").Value = 250

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2467
$ws.Range("F8").Value = 3011
$ws.Range("F10").Value = 250
$ws.Range("F11").Value = 1769
$ws.Range("F12").Value = 110
$ws.Range("F15").Value = 3581
$ws.Range("F17").Value = 1201
$ws.Range("F18").Value = 29
$ws.Range("F21").Value = 21
$ws.Range("F22").Value = 20
$ws.Range("F30").Value = 141
$ws.Range("F31").Value = 141
$ws.Range("F32").Value = 36
$ws.Range("F33").Value = 5
$ws.Range("F34").Value = 194
$ws.Range("F41").Value = 446
$ws.Range("F45").Value = 30
$ws.Range("F50").Value = 940
